$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos values (data refresh from the coinranking feed).
# The Price column (D) holds plain-text cells in the source sheet (t="inlineStr"),
# not numbers, so numeric-looking values are entered with a leading apostrophe to
# force text entry -- matching the original cell type and preserving exact text
# such as trailing zeros ("1.00") that a numeric value would otherwise drop.
$ws.Range("D2").Value = '26.657.19'
$ws.Range("D3").Value = '1.645.25'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'215.89"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = '1.875.28'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.664.80'
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'4.22"
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = "'66.31"
$ws.Range("E16").Value = '  +4.40%  '
$ws.Range("D17").Value = '26.721.86'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '0.0₃0754'
$ws.Range("D19").Value = "'219.54"
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  +10.56%  '
$ws.Range("D25").Value = "'147.32"
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  +2.97%  '
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("D33").Value = "'3.06"
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("D34").Value = '1.287.73'
$ws.Range("E34").Value = '  +6.40%  '
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("E36").Value = '  +6.48%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").Value = "'0.828"
$ws.Range("E38").Value = '  +2.21%  '
$ws.Range("E39").Value = '  +4.31%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = '  +2.13%  '
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").Value = '1.787.22'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = "'93.59"
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("D46").Value = "'60.10"
$ws.Range("E46").Value = '  +9.61%  '
$ws.Range("D47").Value = "'1.62"
$ws.Range("E47").Value = '  +4.15%  '
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("E51").Value = '  -0.58%  '
